$d = $word.ActiveDocument

# 1. Remove the leftover "_GoBack" bookmark (cursor-position marker Word
#    leaves behind after an editing session).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Fill in the name field: the placeholder line
#       "Name: _____________________________ Date Submitted: _____________"
#    becomes
#       "Name: ________Rajesh Patel_________ Date Submitted: _____________"
#    i.e. the student typed their name over some of the blank underscores.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Name: _____________________________ Date Submitted: _____________", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $mid = $d.Range($start + 14, $start + 26)
    $mid.Text = "Rajesh Patel"
}
